# Error Calculations and Plots
# Removes two rows that no longer belong in the missing-data sample (RM 232
# and SC 92), which shifts all subsequent rows up, and refreshes a batch of
# individual cell values (simulating updated error calculations / newly
# imputed vs. newly removed data points).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232"); this shifts "SC 92" (originally row 28) up to row 27.
$ws.Rows.Item(26).Delete()
# Delete the row that is now "SC 92" (row 27 after the first shift).
$ws.Rows.Item(27).Delete()

# Apply the updated cell values (row numbers below are the FINAL row numbers,
# i.e. after the two row deletions above).
$ws.Range("D3").Value = -14.2
$ws.Range("E4").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("F6").Value = 16.43
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("F12").Value = $null
$ws.Range("F14").Value = 17.76
$ws.Range("E17").Value = $null
$ws.Range("F17").Value = 17.78
$ws.Range("E18").Value = $null
$ws.Range("F19").Value = 17.81
$ws.Range("F20").Value = $null
$ws.Range("D21").Value = -14.3
$ws.Range("D23").Value = $null
$ws.Range("F23").Value = $null
$ws.Range("F25").Value = $null
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("D32").Value = -14.7
